$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$newobj = $ws.Shapes.AddChart2(-1, 75)
try {
  Write-Output ("Placement before: " + $newobj.Placement)
} catch { Write-Output ("err get " + $_) }
try {
  $newobj.Placement = 1
  Write-Output "set ok"
} catch { Write-Output ("err set " + $_) }
